# Addition script: adds a delta quantity onto the existing "Bestand" (stock)
# value in column D for a set of article rows on the "Lagerbestand M0129"
# sheet. Also validates that the expected sheet exists before editing it and
# raises a clear error message if the workbook's sheet is not named as
# expected.

$wb = $excel.ActiveWorkbook

$expectedSheetName = "Lagerbestand M0129"
$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq $expectedSheetName) {
        $ws = $sheet
        break
    }
}

if ($ws -eq $null) {
    throw "Sheet '$expectedSheetName' was not found - the workbook's Sheet is not named correctly."
}

# Row (in the sheet) -> quantity to add to the current "Bestand" value in column D.
$additions = [ordered]@{
    2   = -150
    9   = -30
    14  = -180
    15  = 80
    16  = 80
    17  = 80
    18  = 160
    19  = 200
    21  = 100
    22  = -21
    23  = 200
    29  = -204
    30  = -54
    33  = 111
    38  = 200
    40  = -3
    41  = -204
    42  = -204
    43  = -156
    45  = -54
    46  = -6
    48  = 194
    51  = -60
    52  = -18
    53  = -12
    56  = -42
    63  = -15
    67  = 200
    74  = 120
    75  = -45
    76  = -300
    82  = -6
    90  = 245
    91  = 254
    92  = 260
    93  = 251
    97  = -30
    99  = -24
    118 = -12
    124 = -9
    125 = -9
    127 = 157
    128 = 160
    155 = -300
    160 = 251
    180 = 12
    181 = 54
    182 = -6
    193 = 80
    207 = -260
    225 = -6
    230 = -35
    231 = -3
    238 = -300
    242 = -15
    257 = 20
    263 = 40
    266 = 160
}

$col = 4  # column D = "Bestand"

foreach ($row in $additions.Keys) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value2 = $cell.Value2 + $additions[$row]
}
